# Cập nhật Đại Hoàng - Hoa Lư - Bái Dính
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nhap lieu")

# --- Row 38: add new "Mã kho" code (Dầu Điêzen 0,001S Mức 5 column) ---
$ws.Range("H38").Value = "HH050-030"

# --- Row 40: fill in the new "Đại Hoàng" warehouse group ---
# (entry order matches how the values were typed in originally: name, Mã kho, Mã khách, then the rest)
$ws.Range("D40").Value = "Đại Hoàng"
$ws.Range("J40").Value = "KDNL178"
$ws.Range("M40").Value = "KDNL078"
$ws.Range("F40").Value = "HH007044"
$ws.Range("G40").Value = "HH009-001"
$ws.Range("H40").Value = "HH050-034"
$ws.Range("I40").Value = "HH0638"
$ws.Range("K40").Value = "1K25TDH"
$ws.Range("L40").Value = "Nam Định"

# --- Row 41: fill in the new "Hoa Lư" warehouse group ---
$ws.Range("D41").Value = "Hoa Lư"
$ws.Range("J41").Value = "KDNL179"
$ws.Range("M41").Value = "KDNL079"
$ws.Range("F41").Value = "HH007043"
$ws.Range("G41").Value = "HH009-003"
$ws.Range("H41").Value = "HH050-033"
$ws.Range("I41").Value = "HH0639"
$ws.Range("K41").Value = "1K25THL"
$ws.Range("L41").Value = "Nam Định"

# --- Row 42: fill in the new "Bái Đính" warehouse group ---
$ws.Range("D42").Value = "Bái Đính"
$ws.Range("J42").Value = "KDNL177"
$ws.Range("M42").Value = "KDNL077"
$ws.Range("F42").Value = "HH007042"
$ws.Range("G42").Value = "HH009-002"
$ws.Range("H42").Value = "HH050-032"
$ws.Range("I42").Value = "HH06370"
$ws.Range("K42").Value = "1K25TBD"
$ws.Range("L42").Value = "Nam Định"

# Row 40 got a touch shorter once it was filled in
$ws.Rows("40").RowHeight = 9.75

# Column widths grew to fit the new content
$ws.Columns("A").ColumnWidth = 21.67
$ws.Columns("F").ColumnWidth = 15.82

# Freeze the header row and leave the selection on the last edited cell
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K43").Select()
